$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right under
#    the H1 title ("Play Aquaman Slot for Free in 2021" / "Meta description: …").
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Meta description: Read our unbiased review of Aquaman online slot game and play for free in 2021. Enjoy two respin functions, progressive jackpots, and a unique underwater theme.", `
              $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$metaRange = $find.Parent
$metaRange.Expand(4) | Out-Null   # wdParagraph = 4 -> grab the whole paragraph incl. mark
$metaRange.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new bold "Play Aquaman Slot for Free in 2021" paragraph right
#    before the closing italic "image prompt" paragraph (the last paragraph
#    in the body, just ahead of the section break).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertStart = $lastPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$newHeading = "Play Aquaman Slot for Free in 2021"
$xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $newHeading + '</w:t></w:r></w:p></w:body>' + `
  '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlFragment) | Out-Null

# Split the freshly-inserted text from the old last paragraph by turning it
# into its own paragraph.
$breakPos = $insertStart + $newHeading.Length
$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 3) Swap out the old "Create a feature image…" AI-art prompt text for the
#    new meta-description copy, keeping its italic run formatting intact.
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Create a feature image that perfectly captures the adventurous and fun spirit of the Aquaman slot game! The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. Make sure the Maya warrior is engaged in an exciting activity related to the game, such as spinning the reels with Aquaman or finding treasures in Atlantis. Use vibrant and bold colors to make the image standout, and include some of the game symbols like the trident, Mera, or the Aquaman symbol to tie everything together. Let the image showcase the thrill and excitement of this exciting game and make it the perfect visual representation of the Aquaman slot game.", `
               $true, $false, $false, $false, $false, $true, 1, $false, `
               "Read our unbiased review of Aquaman online slot game and play for free in 2021. Enjoy two respin functions, progressive jackpots, and a unique underwater theme.", 2) | Out-Null

Write-Host "Paragraph count:" $d.Paragraphs.Count
